$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Update cell values ---
# Sheet: Overview
$wsOverview.Range("A2").Value = "ffff7ecdeb47-bfab-4196-9c9f-484d19730d84.md"
$wsOverview.Range("A3").Value = "ffffffd90b9004-14d3-4dff-8011-245a1cb20459.md"
$wsOverview.Range("A4").Value = "0c774a1e-471f-4194-b146-314f38ec58bb.md"
$wsOverview.Range("B4").Value = "Ready for handoff"
$wsOverview.Range("C4").Value = "Ready for handoff"

# Sheet: zh-cn
$wsZhCn.Range("A2").Value = "ffff7ecdeb47-bfab-4196-9c9f-484d19730d84.md"
$wsZhCn.Range("C2").Value = "836d1a28-a719-4178-8818-6ed2eabf6198.af2f589f26d8f6ee9fe841d6a9a3da353b03b62c.zh-cn.xlf"
$wsZhCn.Range("D2").Value = "2016-01-25 11:17:52"
$wsZhCn.Range("E2").Value = "836d1a28-a719-4178-8818-6ed2eabf6198.md"
$wsZhCn.Range("F2").Value = "836d1a28-a719-4178-8818-6ed2eabf6198.af2f589f26d8f6ee9fe841d6a9a3da353b03b62c.zh-cn.xlf"
$wsZhCn.Range("G2").Value = "2016-01-25 11:18:34"
$wsZhCn.Range("A3").Value = "ffffffd90b9004-14d3-4dff-8011-245a1cb20459.md"
$wsZhCn.Range("A4").Value = "0c774a1e-471f-4194-b146-314f38ec58bb.md"
$wsZhCn.Range("B4").Value = "Ready for handoff"
$wsZhCn.Range("C4").Value = "0c774a1e-471f-4194-b146-314f38ec58bb.9ba72ff16b4108bd2e09ed1788086f93a297a383.zh-cn.xlf"
$wsZhCn.Range("D4").Value = "2016-01-25 11:21:21"
$wsZhCn.Range("E4").Value = "0c774a1e-471f-4194-b146-314f38ec58bb.md"
$wsZhCn.Range("F4").Value = "0c774a1e-471f-4194-b146-314f38ec58bb.9ba72ff16b4108bd2e09ed1788086f93a297a383.zh-cn.xlf"
$wsZhCn.Range("G4").Value = "2016-01-25 11:20:28"

# Sheet: de-de
$wsDeDe.Range("A2").Value = "ffff7ecdeb47-bfab-4196-9c9f-484d19730d84.md"
$wsDeDe.Range("C2").Value = "836d1a28-a719-4178-8818-6ed2eabf6198.af2f589f26d8f6ee9fe841d6a9a3da353b03b62c.de-de.xlf"
$wsDeDe.Range("D2").Value = "2016-01-25 11:18:03"
$wsDeDe.Range("E2").Value = "836d1a28-a719-4178-8818-6ed2eabf6198.md"
$wsDeDe.Range("F2").Value = "836d1a28-a719-4178-8818-6ed2eabf6198.af2f589f26d8f6ee9fe841d6a9a3da353b03b62c.de-de.xlf"
$wsDeDe.Range("G2").Value = "2016-01-25 11:18:50"
$wsDeDe.Range("A3").Value = "ffffffd90b9004-14d3-4dff-8011-245a1cb20459.md"
$wsDeDe.Range("A4").Value = "0c774a1e-471f-4194-b146-314f38ec58bb.md"
$wsDeDe.Range("B4").Value = "Ready for handoff"
$wsDeDe.Range("C4").Value = "0c774a1e-471f-4194-b146-314f38ec58bb.9ba72ff16b4108bd2e09ed1788086f93a297a383.de-de.xlf"
$wsDeDe.Range("D4").Value = "2016-01-25 11:21:30"
$wsDeDe.Range("E4").Value = "0c774a1e-471f-4194-b146-314f38ec58bb.md"
$wsDeDe.Range("F4").Value = "0c774a1e-471f-4194-b146-314f38ec58bb.9ba72ff16b4108bd2e09ed1788086f93a297a383.de-de.xlf"
$wsDeDe.Range("G4").Value = "2016-01-25 11:20:44"

# --- Rebuild hyperlinks (preserve same target URLs / rIds, update display text) ---
# Sheet: Overview
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e3304c3d9c196f781205d7896c7f198ab95a5db8/e2e/0c774a1e-471f-4194-b146-314f38ec58bb.md", "", "", "ffff7ecdeb47-bfab-4196-9c9f-484d19730d84.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e3304c3d9c196f781205d7896c7f198ab95a5db8/e2e/ffff7ecdeb47-bfab-4196-9c9f-484d19730d84.md", "", "", "ffffffd90b9004-14d3-4dff-8011-245a1cb20459.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/e3304c3d9c196f781205d7896c7f198ab95a5db8/e2e/ffffffd90b9004-14d3-4dff-8011-245a1cb20459.md", "", "", "0c774a1e-471f-4194-b146-314f38ec58bb.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/e3304c3d9c196f781205d7896c7f198ab95a5db8/.localization-config", "", "", ".localization-config") | Out-Null

# Sheet: zh-cn
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e3304c3d9c196f781205d7896c7f198ab95a5db8/e2e/0c774a1e-471f-4194-b146-314f38ec58bb.md", "", "", "ffff7ecdeb47-bfab-4196-9c9f-484d19730d84.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/da65cedf21253fa50ac608762984cfe1df703c1e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/0c774a1e-471f-4194-b146-314f38ec58bb.9ba72ff16b4108bd2e09ed1788086f93a297a383.zh-cn.xlf", "", "", "836d1a28-a719-4178-8818-6ed2eabf6198.af2f589f26d8f6ee9fe841d6a9a3da353b03b62c.zh-cn.xlf") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/5f24c0aaf64b77197cd2335dae3ea36510a3deb5/e2e/0c774a1e-471f-4194-b146-314f38ec58bb.md", "", "", "836d1a28-a719-4178-8818-6ed2eabf6198.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/cfb53e96f34acf816626e721e509fbf11747352e/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/0c774a1e-471f-4194-b146-314f38ec58bb.9ba72ff16b4108bd2e09ed1788086f93a297a383.zh-cn.xlf", "", "", "836d1a28-a719-4178-8818-6ed2eabf6198.af2f589f26d8f6ee9fe841d6a9a3da353b03b62c.zh-cn.xlf") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e3304c3d9c196f781205d7896c7f198ab95a5db8/e2e/ffff7ecdeb47-bfab-4196-9c9f-484d19730d84.md", "", "", "ffffffd90b9004-14d3-4dff-8011-245a1cb20459.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/edeca7c8cedf741e2ec4f04014835c4009a1799e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/836d1a28-a719-4178-8818-6ed2eabf6198.af2f589f26d8f6ee9fe841d6a9a3da353b03b62c.zh-cn.xlf", "", "", "836d1a28-a719-4178-8818-6ed2eabf6198.af2f589f26d8f6ee9fe841d6a9a3da353b03b62c.zh-cn.xlf") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/e33d83719dd61434dab207f6b624497457b7aa97/e2e/836d1a28-a719-4178-8818-6ed2eabf6198.md", "", "", "836d1a28-a719-4178-8818-6ed2eabf6198.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/7d63ee76760b7873138a1e3af35c199f47a3cb66/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/836d1a28-a719-4178-8818-6ed2eabf6198.af2f589f26d8f6ee9fe841d6a9a3da353b03b62c.zh-cn.xlf", "", "", "836d1a28-a719-4178-8818-6ed2eabf6198.af2f589f26d8f6ee9fe841d6a9a3da353b03b62c.zh-cn.xlf") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/e3304c3d9c196f781205d7896c7f198ab95a5db8/e2e/ffffffd90b9004-14d3-4dff-8011-245a1cb20459.md", "", "", "0c774a1e-471f-4194-b146-314f38ec58bb.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/edeca7c8cedf741e2ec4f04014835c4009a1799e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/836d1a28-a719-4178-8818-6ed2eabf6198.af2f589f26d8f6ee9fe841d6a9a3da353b03b62c.zh-cn.xlf", "", "", "0c774a1e-471f-4194-b146-314f38ec58bb.9ba72ff16b4108bd2e09ed1788086f93a297a383.zh-cn.xlf") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/e33d83719dd61434dab207f6b624497457b7aa97/e2e/836d1a28-a719-4178-8818-6ed2eabf6198.md", "", "", "0c774a1e-471f-4194-b146-314f38ec58bb.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/7d63ee76760b7873138a1e3af35c199f47a3cb66/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/836d1a28-a719-4178-8818-6ed2eabf6198.af2f589f26d8f6ee9fe841d6a9a3da353b03b62c.zh-cn.xlf", "", "", "0c774a1e-471f-4194-b146-314f38ec58bb.9ba72ff16b4108bd2e09ed1788086f93a297a383.zh-cn.xlf") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/e3304c3d9c196f781205d7896c7f198ab95a5db8/.localization-config", "", "", ".localization-config") | Out-Null

# Sheet: de-de
$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e3304c3d9c196f781205d7896c7f198ab95a5db8/e2e/0c774a1e-471f-4194-b146-314f38ec58bb.md", "", "", "ffff7ecdeb47-bfab-4196-9c9f-484d19730d84.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c2d4f1fd73c669ba14ade9e9d7c1df035de90dea/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/0c774a1e-471f-4194-b146-314f38ec58bb.9ba72ff16b4108bd2e09ed1788086f93a297a383.de-de.xlf", "", "", "836d1a28-a719-4178-8818-6ed2eabf6198.af2f589f26d8f6ee9fe841d6a9a3da353b03b62c.de-de.xlf") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/4d826a379150669d3d65ed394afa52768aa9750a/e2e/0c774a1e-471f-4194-b146-314f38ec58bb.md", "", "", "836d1a28-a719-4178-8818-6ed2eabf6198.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/57784fe9ed492bf66b6108e0fb61ca56133ad962/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/0c774a1e-471f-4194-b146-314f38ec58bb.9ba72ff16b4108bd2e09ed1788086f93a297a383.de-de.xlf", "", "", "836d1a28-a719-4178-8818-6ed2eabf6198.af2f589f26d8f6ee9fe841d6a9a3da353b03b62c.de-de.xlf") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e3304c3d9c196f781205d7896c7f198ab95a5db8/e2e/ffff7ecdeb47-bfab-4196-9c9f-484d19730d84.md", "", "", "ffffffd90b9004-14d3-4dff-8011-245a1cb20459.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bf411da73484f28ca38449b8fe92e9eb2e475988/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/836d1a28-a719-4178-8818-6ed2eabf6198.af2f589f26d8f6ee9fe841d6a9a3da353b03b62c.de-de.xlf", "", "", "836d1a28-a719-4178-8818-6ed2eabf6198.af2f589f26d8f6ee9fe841d6a9a3da353b03b62c.de-de.xlf") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/039698ffc4da609c81b74922efe1b957177065c7/e2e/836d1a28-a719-4178-8818-6ed2eabf6198.md", "", "", "836d1a28-a719-4178-8818-6ed2eabf6198.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e8f1bd4024e2a4af7f0b0662d86311d3a20313a5/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/836d1a28-a719-4178-8818-6ed2eabf6198.af2f589f26d8f6ee9fe841d6a9a3da353b03b62c.de-de.xlf", "", "", "836d1a28-a719-4178-8818-6ed2eabf6198.af2f589f26d8f6ee9fe841d6a9a3da353b03b62c.de-de.xlf") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/e3304c3d9c196f781205d7896c7f198ab95a5db8/e2e/ffffffd90b9004-14d3-4dff-8011-245a1cb20459.md", "", "", "0c774a1e-471f-4194-b146-314f38ec58bb.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bf411da73484f28ca38449b8fe92e9eb2e475988/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/836d1a28-a719-4178-8818-6ed2eabf6198.af2f589f26d8f6ee9fe841d6a9a3da353b03b62c.de-de.xlf", "", "", "0c774a1e-471f-4194-b146-314f38ec58bb.9ba72ff16b4108bd2e09ed1788086f93a297a383.de-de.xlf") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/039698ffc4da609c81b74922efe1b957177065c7/e2e/836d1a28-a719-4178-8818-6ed2eabf6198.md", "", "", "0c774a1e-471f-4194-b146-314f38ec58bb.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e8f1bd4024e2a4af7f0b0662d86311d3a20313a5/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/836d1a28-a719-4178-8818-6ed2eabf6198.af2f589f26d8f6ee9fe841d6a9a3da353b03b62c.de-de.xlf", "", "", "0c774a1e-471f-4194-b146-314f38ec58bb.9ba72ff16b4108bd2e09ed1788086f93a297a383.de-de.xlf") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/e3304c3d9c196f781205d7896c7f198ab95a5db8/.localization-config", "", "", ".localization-config") | Out-Null
